# On-call tracker update: mark an additional on-call entry (tally of 1)
# in column C - the first day of each period - for several teachers
# who picked up an extra on-call shift. This adds new data cells that
# did not exist before, and the sheet's existing SUM()/MAX() formulas
# (columns X, Y, Z, AA and the period/grand totals in rows 71-76)
# recalculate automatically to reflect the new counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellsToFill = @("C6", "C8", "C10", "C13", "C15", "C43", "C44", "C45", "C46")

foreach ($cellRef in $cellsToFill) {
    $ws.Range($cellRef).Value = 1
}
